$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "66.231.96"
Set-TextValue "E2" "  -0.08%  "
Set-TextValue "D3" "3.566.08"
Set-TextValue "E3" "  +1.28%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "606.37"
Set-TextValue "E6" "  -0.72%  "
Set-TextValue "D7" "3.565.75"
Set-TextValue "E7" "  +1.27%  "
Set-TextValue "E8" "  +0.08%  "
Set-TextValue "E9" "  +2.51%  "
Set-TextValue "E10" "  -0.04%  "
Set-TextValue "D11" "7.81"
Set-TextValue "E11" "  -2.50%  "
Set-TextValue "E12" "  +0.03%  "
Set-TextValue "D13" "4.169.75"
Set-TextValue "E13" "  +1.28%  "
Set-TextValue "E14" "  +0.13%  "
Set-TextValue "D15" "30.39"
Set-TextValue "E15" "  -0.04%  "
Set-TextValue "D16" "3.573.39"
Set-TextValue "E16" "  +1.59%  "
Set-TextValue "D17" "66.291.95"
Set-TextValue "E17" "  -0.04%  "
Set-TextValue "E18" "  -0.64%  "
Set-TextValue "D19" "11.39"
Set-TextValue "E19" "  +7.06%  "
Set-TextValue "D20" "6.22"
Set-TextValue "E20" "  +0.27%  "
Set-TextValue "D21" "14.77"
Set-TextValue "E21" "  -1.09%  "
Set-TextValue "D22" "431.21"
Set-TextValue "E22" "  +1.17%  "
Set-TextValue "D23" "0.613"
Set-TextValue "E23" "  +2.29%  "
Set-TextValue "D24" "79.48"
Set-TextValue "E24" "  +1.70%  "
Set-TextValue "D25" "3.708.03"
Set-TextValue "E25" "  +1.42%  "
Set-TextValue "E26" "  +0.02%  "
Set-TextValue "E27" "  -1.91%  "
Set-TextValue "E28" "  +0.97%  "
Set-TextValue "D29" "9.17"
Set-TextValue "E29" "  -1.23%  "
Set-TextValue "D30" "7.92"
Set-TextValue "E30" "  -1.39%  "
Set-TextValue "D31" "0.986"
Set-TextValue "E31" "  -1.42%  "
Set-TextValue "D32" "3.561.93"
Set-TextValue "E32" "  +1.68%  "
Set-TextValue "D33" "25.46"
Set-TextValue "E33" "  +0.82%  "
Set-TextValue "E35" "  -7.99%  "
Set-TextValue "D36" "7.86"
Set-TextValue "E38" "  -0.69%  "
Set-TextValue "D39" "5.61"
Set-TextValue "E39" "  +0.02%  "
Set-TextValue "D40" "172.87"
Set-TextValue "E40" "  +1.24%  "
Set-TextValue "E41" "  -0.82%  "
Set-TextValue "E42" "  +0.36%  "
Set-TextValue "D43" "0.889"
Set-TextValue "E43" "  +0.05%  "
Set-TextValue "E44" "  +2.08%  "
Set-TextValue "D45" "45.99"
Set-TextValue "E45" "  +1.13%  "
Set-TextValue "D46" "1.00"
Set-TextValue "E46" "  +0.06%  "
Set-TextValue "E47" "  -1.44%  "
Set-TextValue "E48" "  +1.11%  "
Set-TextValue "D49" "24.88"
Set-TextValue "E49" "  -4.13%  "
Set-TextValue "D51" "23.27"
Set-TextValue "E51" "  +4.50%  "
